$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 11.98327633333333
$ws.Range("H2").Value = 35.949829
$ws.Range("I2").Value = 0.03345300399843466
$ws.Range("J2").Value = 0.03345300399843466
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 153.5290173333333
$ws.Range("N2").Value = 460.587052
$ws.Range("O2").Value = 0.3172206968818489
$ws.Range("P2").Value = 0.317220696881849
$ws.Range("Q2").Value = 1839.780639890456
$ws.Range("R2").Value = 16558.02575901411
$ws.Range("S2").Value = 0.01061198524117472
$ws.Range("T2").Value = 0.01061198524117472
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 11.98327633333333
$ws.Range("H3").Value = 35.949829
$ws.Range("I3").Value = 0.03345300399843466
$ws.Range("J3").Value = 0.03345300399843466
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 168.7997026666667
$ws.Range("N3").Value = 506.3991080000001
$ws.Range("O3").Value = 0.3487728915577651
$ws.Range("P3").Value = 0.3487728915577651
$ws.Range("Q3").Value = 2022.773482039171
$ws.Range("R3").Value = 18204.96133835254
$ws.Range("S3").Value = 0.01166750093582753
$ws.Range("T3").Value = 0.01166750093582753
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 11.98327633333333
$ws.Range("H4").Value = 35.949829
$ws.Range("I4").Value = 0.03345300399843466
$ws.Range("J4").Value = 0.03345300399843466
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 68.09032333333333
$ws.Range("N4").Value = 204.27097
$ws.Range("O4").Value = 0.1406878008722904
$ws.Range("P4").Value = 0.1406878008722904
$ws.Range("Q4").Value = 815.9451601293478
$ws.Range("R4").Value = 7343.506441164131
$ws.Range("S4").Value = 0.004706429565111709
$ws.Range("T4").Value = 0.00470642956511171
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 11.98327633333333
$ws.Range("H5").Value = 35.949829
$ws.Range("I5").Value = 0.03345300399843466
$ws.Range("J5").Value = 0.03345300399843466
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 93.562673
$ws.Range("N5").Value = 280.688019
$ws.Range("O5").Value = 0.1933186106880956
$ws.Range("P5").Value = 0.1933186106880956
$ws.Range("Q5").Value = 1121.187365044306
$ws.Range("R5").Value = 10090.68628539875
$ws.Range("S5").Value = 0.006467088256320696
$ws.Range("T5").Value = 0.006467088256320696
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 332.3726806666667
$ws.Range("H6").Value = 997.1180420000001
$ws.Range("I6").Value = 0.9278651602470024
$ws.Range("J6").Value = 0.9278651602470025
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 153.5290173333333
$ws.Range("N6").Value = 460.587052
$ws.Range("O6").Value = 0.3172206968818489
$ws.Range("P6").Value = 0.317220696881849
$ws.Range("Q6").Value = 51028.85105119913
$ws.Range("R6").Value = 459259.6594607922
$ws.Range("S6").Value = 0.2943380327459425
$ws.Range("T6").Value = 0.2943380327459426
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 332.3726806666667
$ws.Range("H7").Value = 997.1180420000001
$ws.Range("I7").Value = 0.9278651602470024
$ws.Range("J7").Value = 0.9278651602470025
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 168.7997026666667
$ws.Range("N7").Value = 506.3991080000001
$ws.Range("O7").Value = 0.3487728915577651
$ws.Range("P7").Value = 0.3487728915577651
$ws.Range("Q7").Value = 56104.40967105629
$ws.Range("R7").Value = 504939.6870395066
$ws.Range("S7").Value = 0.3236142149150561
$ws.Range("T7").Value = 0.3236142149150561
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 332.3726806666667
$ws.Range("H8").Value = 997.1180420000001
$ws.Range("I8").Value = 0.9278651602470024
$ws.Range("J8").Value = 0.9278651602470025
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 68.09032333333333
$ws.Range("N8").Value = 204.27097
$ws.Range("O8").Value = 0.1406878008722904
$ws.Range("P8").Value = 0.1406878008722904
$ws.Range("Q8").Value = 22631.36329376008
$ws.Range("R8").Value = 203682.2696438408
$ws.Range("S8").Value = 0.1305393089011661
$ws.Range("T8").Value = 0.1305393089011661
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 332.3726806666667
$ws.Range("H9").Value = 997.1180420000001
$ws.Range("I9").Value = 0.9278651602470024
$ws.Range("J9").Value = 0.9278651602470025
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 93.562673
$ws.Range("N9").Value = 280.688019
$ws.Range("O9").Value = 0.1933186106880956
$ws.Range("P9").Value = 0.1933186106880956
$ws.Range("Q9").Value = 31097.67643534876
$ws.Range("R9").Value = 279879.0879181388
$ws.Range("S9").Value = 0.1793736036848377
$ws.Range("T9").Value = 0.1793736036848377
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.08615933333333332
$ws.Range("H10").Value = 0.258478
$ws.Range("I10").Value = 0.0002405259164795302
$ws.Range("J10").Value = 0.0002405259164795302
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 153.5290173333333
$ws.Range("N10").Value = 460.587052
$ws.Range("O10").Value = 0.3172206968818489
$ws.Range("P10").Value = 0.317220696881849
$ws.Range("Q10").Value = 13.22795778076178
$ws.Range("R10").Value = 119.051620026856
$ws.Range("S10").Value = 0.00007629979884378195
$ws.Range("T10").Value = 0.00007629979884378196
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.08615933333333332
$ws.Range("H11").Value = 0.258478
$ws.Range("I11").Value = 0.0002405259164795302
$ws.Range("J11").Value = 0.0002405259164795302
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 168.7997026666667
$ws.Range("N11").Value = 506.3991080000001
$ws.Range("O11").Value = 0.3487728915577651
$ws.Range("P11").Value = 0.3487728915577651
$ws.Range("Q11").Value = 14.54366984862489
$ws.Range("R11").Value = 130.893028637624
$ws.Range("S11").Value = 0.00008388891938514723
$ws.Range("T11").Value = 0.00008388891938514724
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.08615933333333332
$ws.Range("H12").Value = 0.258478
$ws.Range("I12").Value = 0.0002405259164795302
$ws.Range("J12").Value = 0.0002405259164795302
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 68.09032333333333
$ws.Range("N12").Value = 204.27097
$ws.Range("O12").Value = 0.1406878008722904
$ws.Range("P12").Value = 0.1406878008722904
$ws.Range("Q12").Value = 5.86661686485111
$ws.Range("R12").Value = 52.79955178366
$ws.Range("S12").Value = 0.00003383906224229729
$ws.Range("T12").Value = 0.00003383906224229729
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.08615933333333332
$ws.Range("H13").Value = 0.258478
$ws.Range("I13").Value = 0.0002405259164795302
$ws.Range("J13").Value = 0.0002405259164795302
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 93.562673
$ws.Range("N13").Value = 280.688019
$ws.Range("O13").Value = 0.1933186106880956
$ws.Range("P13").Value = 0.1933186106880956
$ws.Range("Q13").Value = 8.061297530564666
$ws.Range("R13").Value = 72.551677775082
$ws.Range("S13").Value = 0.00004649813600830369
$ws.Range("T13").Value = 0.0000464981360083037
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 13.77014866666667
$ws.Range("H14").Value = 41.310446
$ws.Range("I14").Value = 0.03844130983808348
$ws.Range("J14").Value = 0.03844130983808348
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 153.5290173333333
$ws.Range("N14").Value = 460.587052
$ws.Range("O14").Value = 0.3172206968818489
$ws.Range("P14").Value = 0.317220696881849
$ws.Range("Q14").Value = 2114.117393327243
$ws.Range("R14").Value = 19027.05653994519
$ws.Range("S14").Value = 0.01219437909588791
$ws.Range("T14").Value = 0.01219437909588792
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 13.77014866666667
$ws.Range("H15").Value = 41.310446
$ws.Range("I15").Value = 0.03844130983808348
$ws.Range("J15").Value = 0.03844130983808348
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 168.7997026666667
$ws.Range("N15").Value = 506.3991080000001
$ws.Range("O15").Value = 0.3487728915577651
$ws.Range("P15").Value = 0.3487728915577651
$ws.Range("Q15").Value = 2324.39700060913
$ws.Range("R15").Value = 20919.57300548217
$ws.Range("S15").Value = 0.01340728678749634
$ws.Range("T15").Value = 0.01340728678749634
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 13.77014866666667
$ws.Range("H16").Value = 41.310446
$ws.Range("I16").Value = 0.03844130983808348
$ws.Range("J16").Value = 0.03844130983808348
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 68.09032333333333
$ws.Range("N16").Value = 204.27097
$ws.Range("O16").Value = 0.1406878008722904
$ws.Range("P16").Value = 0.1406878008722904
$ws.Range("Q16").Value = 937.6138750614022
$ws.Range("R16").Value = 8438.524875552621
$ws.Range("S16").Value = 0.005408223343770306
$ws.Range("T16").Value = 0.005408223343770308
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 13.77014866666667
$ws.Range("H17").Value = 41.310446
$ws.Range("I17").Value = 0.03844130983808348
$ws.Range("J17").Value = 0.03844130983808348
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 93.562673
$ws.Range("N17").Value = 280.688019
$ws.Range("O17").Value = 0.1933186106880956
$ws.Range("P17").Value = 0.1933186106880956
$ws.Range("Q17").Value = 1288.371916860719
$ws.Range("R17").Value = 11595.34725174647
$ws.Range("S17").Value = 0.007431420610928921
$ws.Range("T17").Value = 0.007431420610928922

Write-Output "updated cells"
